$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 850
$ws.Range("B3").Value = 100
$ws.Range("B4").Value = 75
$ws.Range("B5").Value = 75
$ws.Range("B6").Value = 220
